# Update cryptos list (price/volume figures) and swap TRON/Chainlink row order
# as scraped on Sat May 25 17:42:53 UTC 2024.
#
# Numeric-looking price values are prefixed with a leading apostrophe so
# Excel stores them verbatim as text (preserving trailing zeros / exact
# digits) instead of auto-converting them to floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.010.29"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "3.741.43"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'601.90"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'167.15"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "3.740.72"
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("E10").Value = "  +2.48%  "
$ws.Range("E11").Value = "  +2.99%  "
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "'37.94"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("E14").Value = "  +1.80%  "
$ws.Range("D15").Value = "4.367.22"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "3.745.74"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "69.006.76"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.113"
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'17.26"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "'10.99"
$ws.Range("E21").Value = "  +19.62%  "
$ws.Range("D22").Value = "'492.16"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "'0.725"
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("D24").Value = "'0.0000152"
$ws.Range("E24").Value = "  +8.25%  "
$ws.Range("D25").Value = "'84.78"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").Value = "'2.30"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "'12.33"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("D28").Value = "'10.09"
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +2.19%  "
$ws.Range("D31").Value = "'2.48"
$ws.Range("E31").Value = "  +4.23%  "
$ws.Range("D32").Value = "'8.04"
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("D33").Value = "'31.48"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "3.886.96"
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").Value = "3.677.84"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("E39").Value = "  +3.14%  "
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("E42").Value = "  +5.83%  "
$ws.Range("D43").Value = "'431.56"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "'48.64"
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("D46").Value = "'8.47"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D48").Value = "'40.28"
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("D49").Value = "'141.17"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").Value = "2.783.45"
$ws.Range("E50").Value = "  +1.71%  "
$ws.Range("E51").Value = "  +0.77%  "
